$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Text
    $cell.Value = "V-" + $old
}
